# "ndimas profile provider, main page"
# Add latitude/longitude columns (E, F) to each address row, drop the
# leftover fill style on A3, move the selection, and set the page to
# portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New E/F data columns
$ws.Range("E1").Value = 3.4996999999999998
$ws.Range("F1").Value = 98.9923

$ws.Range("E2").Value = 3.5861999999999998
$ws.Range("F2").Value = 98.600300000000004

$ws.Range("E3").Value = 4.1567999999999996
$ws.Range("F3").Value = 98.600499999999997

# A3 no longer carries the (unused) fill style
$ws.Range("A3").ClearFormats()

# Page orientation -> portrait
$ws.PageSetup.Orientation = 1

# Scroll the view over and land the selection on H3
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$ws.Range("H3").Select()
